$d = $word.ActiveDocument

# 1) Collapse the "{{ case_number }}" placeholder (CASE NO. line) to a single run.
$d.Content.Find.Execute("{{ case_number }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ case_number }}", 2)

# 2) Collapse the "{{ defendant_name }}" placeholder wherever it appears in the body.
$d.Content.Find.Execute("{{ defendant_name }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ defendant_name }}", 2)

# 3) Collapse "{{ plea_trial_date }}" placeholder.
$d.Content.Find.Execute("{{ plea_trial_date }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ plea_trial_date }}", 2)

# 4) Collapse the "{%tc for charge in charges_list %}" / "{%tc endfor %}" / "{{ charge.xxx }}" runs.
$d.Content.Find.Execute("{%tc for charge in charges_list %}", $true, $false, $false, $false, $false, $true, 1, $false, "{%tc for charge in charges_list %}", 2)
$d.Content.Find.Execute("{%tc endfor %}", $true, $false, $false, $false, $false, $true, 1, $false, "{%tc endfor %}", 2)
$d.Content.Find.Execute("{{ charge.offense }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ charge.offense }}", 2)
$d.Content.Find.Execute("{{ charge.degree }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ charge.degree }}", 2)
$d.Content.Find.Execute("{{ charge.plea}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ charge.plea}}", 2)
$d.Content.Find.Execute("{{ charge.finding }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ charge.finding }}", 2)
$d.Content.Find.Execute("{{ charge.fines_amount }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ charge.fines_amount }}", 2)
$d.Content.Find.Execute("{{ charge.fines_suspended }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ charge.fines_suspended }}", 2)
$d.Content.Find.Execute("{{ charge.court_costs }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ charge.court_costs }}", 2)

# 5) Expand the ability-to-pay sentence with the new balance-due clause.
$d.Content.Find.Execute("{{ ability_to_pay_time }}. Community", $true, $false, $false, $false, $false, $true, 1, $false, "{{ ability_to_pay_time }} and absent further court order the total balance of fines and costs shall be paid by {{ balance_due_date }}. Community", 2)

# 6) Shorten "ADMINISTRATIVE JUDGE" signature line to "JUDGE".
$d.Content.Find.Execute("ADMINISTRATIVE JUDGE", $true, $false, $false, $false, $false, $true, 1, $false, "JUDGE", 2)

# 7) Footer: collapse the "{{ case_number }}" placeholder runs.
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$ftr.Range.Find.Execute("{{ case_number }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ case_number }}", 2)
